# Fruta / hortaliza, semanal
# Insert a new weekly record as row 211 (pushing the existing rows 211-239
# down to 212-240), matching the new latest "Vega Modelo de Temuco - Mango"
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 211; Excel shifts rows 211:239 down to 212:240
# and carries row formatting (e.g. the date style on column D) along.
$ws.Rows(211).Insert()

# Populate the newly inserted row 211 with the new data record.
$ws.Range("A211").Value = 10
$ws.Range("B211").Value = "Vega Modelo de Temuco"
$ws.Range("C211").Value = "La Araucanía"
$ws.Range("D211").Value = 44505
$ws.Range("E211").Value = 9
$ws.Range("F211").Value = "Fruta"
$ws.Range("G211").Value = 100108
$ws.Range("H211").Value = "Tropicales y subtropicales"
$ws.Range("I211").Value = 100108002
$ws.Range("J211").Value = "Mango"
$ws.Range("K211").Value = "Sin especificar"
$ws.Range("L211").Value = "Primera"
$ws.Range("M211").Value = 155
$ws.Range("N211").Value = 8000
$ws.Range("O211").Value = 8000
$ws.Range("P211").Value = 8000
$ws.Range("Q211").Value = "`$/bandeja 4 kilos"
$ws.Range("R211").Value = "Perú"
$ws.Range("S211").Value = 2000
$ws.Range("T211").Value = 4
